# Reorganize the "道具" sheet:
#   - move the header row (编号/名称/类型) from row 3 up to row 1
#   - the two placeholder rows that used to be rows 1-2 ("1 1 1" / "1 0 1")
#     shift down to rows 2-3
#   - row 4 (id/name/type) is unchanged
#   - row 5's C-style type names are simplified (uint32_t->uint,
#     std::string->string, uint8_t->int)
#   - row 6 (1/钻石/1) is unchanged
#   - final selection ends on C5 (last cell touched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("道具")

$ws.Range("A1").Value = "编号"
$ws.Range("B1").Value = "名称"
$ws.Range("C1").Value = "类型"

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "1"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "0"
$ws.Range("C3").Value = "1"

$ws.Range("A5").Value = "uint"
$ws.Range("B5").Value = "string"
$ws.Range("C5").Value = "int"

$ws.Range("C5").Select()
